$d = $word.ActiveDocument

# Locate the unique phrase "Install required packages". In the source
# document it lives in a single run together with a leading <w:tab/>:
#   <w:r><w:rPr>...</w:rPr><w:tab/><w:t>Install required packages</w:t></w:r>
# The target edit splits that run into three runs - a tab-only run, a
# "Load" run, and a " required packages" run - while leaving the rest of
# the paragraph (" and" / " set working directory") untouched.
$found = $d.Content
$found.Find.ClearFormatting()
$ok = $found.Find.Execute("Install required packages", $true, $false, $false,
                           $false, $false, $true, 1, $false, "", 0)
if (-not $ok) {
    throw "Could not find 'Install required packages' in the document."
}

$para = $found.Paragraphs(1)
$paraEnd = $para.Range.End - 1

# Replace from the tab character (one position before "Install") through
# to the end of the paragraph. InsertXML always splices its replacement
# runs in at the tail of the range it is given, so the range must reach
# the end of the paragraph for the new runs to land in the right order;
# the two trailing runs that come after "packages" are therefore
# reproduced verbatim in the replacement payload below.
$rng = $d.Range($found.Start - 1, $paraEnd)

$rPr = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr>'

$runsXml = ('<w:r>' + $rPr + '<w:tab/></w:r>') + `
           ('<w:r>' + $rPr + '<w:t>Load</w:t></w:r>') + `
           ('<w:r>' + $rPr + '<w:t xml:space="preserve"> required packages</w:t></w:r>') + `
           ('<w:r w:rsidR="00115613">' + $rPr + '<w:t xml:space="preserve"> and</w:t></w:r>') + `
           ('<w:r w:rsidR="002A26B0">' + $rPr + '<w:t xml:space="preserve"> set working directory</w:t></w:r>')

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
       '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
       '<pkg:part pkg:name="/part" pkg:contentType="application/xml"><pkg:xmlData>' + `
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
       '<w:body><w:p>' + $runsXml + '</w:p></w:body></w:document>' + `
       '</pkg:xmlData></pkg:part></pkg:package>'

$rng.InsertXML($xml)
